# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method
#
# Appends one more game's worth of per-play/per-drive numbers to the
# long space-separated logs on the YDS and ST sheets, and updates the
# season-total cells on OFF, DEF, ST, TURNS and PEN accordingly.
#
# Note: reading back via ".Value" on this host returns the property's
# reflection descriptor rather than the cell's contents, so ".Value2"
# is used for both reads and writes below.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet: append new values to the long space-separated number lists
# ---------------------------------------------------------------------
$wsYDS = $wb.Worksheets.Item("YDS")

$wsYDS.Range("B2").Value2 = $wsYDS.Range("B2").Value2 + " 0 6 6 -3 7 4 3 0 3 -3 22 3 6 4 2 0 1 14 13 5 4 2 8 0 1 7 4 2 3"
$wsYDS.Range("C2").Value2 = $wsYDS.Range("C2").Value2 + " 3 3 11 3 6 1 1 3 4 4 5 8 1 4 3 8 0 4 10 1 8 3"
$wsYDS.Range("B3").Value2 = $wsYDS.Range("B3").Value2 + " 30 5 12 14 5 19 18 5 8 8 18 13 -3 27 12 13 10"
$wsYDS.Range("C3").Value2 = $wsYDS.Range("C3").Value2 + " 16 10 24 8 25 12 2 12 12 11 5 6 18 9 4 11 1 -4 2 8 10 4 5 9 12 5 10"

# ---------------------------------------------------------------------
# OFF sheet: update season totals (row 2 = Home, row 3 = Road)
# ---------------------------------------------------------------------
$wsOFF = $wb.Worksheets.Item("OFF")

$wsOFF.Range("C2").Value2 = 216
$wsOFF.Range("D2").Value2 = 15
$wsOFF.Range("E2").Value2 = 11
$wsOFF.Range("F2").Value2 = 73
$wsOFF.Range("G2").Value2 = 73
$wsOFF.Range("I2").Value2 = 14
$wsOFF.Range("J2").Value2 = 44
$wsOFF.Range("N2").Value2 = 12
$wsOFF.Range("O2").Value2 = 22
$wsOFF.Range("P2").Value2 = 9

$wsOFF.Range("B3").Value2 = 14
$wsOFF.Range("C3").Value2 = 140
$wsOFF.Range("F3").Value2 = 91
$wsOFF.Range("G3").Value2 = 23
$wsOFF.Range("H3").Value2 = 22
$wsOFF.Range("I3").Value2 = 36
$wsOFF.Range("J3").Value2 = 58
$wsOFF.Range("L3").Value2 = 262
$wsOFF.Range("M3").Value2 = 168
$wsOFF.Range("Q3").Value2 = 550

# ---------------------------------------------------------------------
# DEF sheet: update season totals (row 2 = Home, row 3 = Road)
# ---------------------------------------------------------------------
$wsDEF = $wb.Worksheets.Item("DEF")

$wsDEF.Range("C2").Value2 = 192
$wsDEF.Range("D2").Value2 = 13
$wsDEF.Range("F2").Value2 = 55
$wsDEF.Range("G2").Value2 = 62
$wsDEF.Range("J2").Value2 = 30
$wsDEF.Range("N2").Value2 = 22
$wsDEF.Range("O2").Value2 = 26

$wsDEF.Range("B3").Value2 = 12
$wsDEF.Range("C3").Value2 = 187
$wsDEF.Range("E3").Value2 = 31
$wsDEF.Range("F3").Value2 = 115
$wsDEF.Range("H3").Value2 = 27
$wsDEF.Range("I3").Value2 = 62
$wsDEF.Range("J3").Value2 = 54
$wsDEF.Range("L3").Value2 = 304
$wsDEF.Range("M3").Value2 = 210
$wsDEF.Range("Q3").Value2 = 549

# ---------------------------------------------------------------------
# ST sheet: append to the long per-kick logs + update season totals
# ---------------------------------------------------------------------
$wsST = $wb.Worksheets.Item("ST")

$wsST.Range("B4").Value2 = $wsST.Range("B4").Value2 + " 66 67"
$wsST.Range("B5").Value2 = $wsST.Range("B5").Value2 + " 34 23"
$wsST.Range("B6").Value2 = $wsST.Range("B6").Value2 + " 21 0"
$wsST.Range("D3").Value2 = $wsST.Range("D3").Value2 + " 39 35"
$wsST.Range("D4").Value2 = $wsST.Range("D4").Value2 + " 0 0"
$wsST.Range("D5").Value2 = $wsST.Range("D5").Value2 + " 0 7"

$wsST.Range("B2").Value2 = 87
$wsST.Range("D2").Value2 = 51
$wsST.Range("F2").Value2 = 182
$wsST.Range("G2").Value2 = 173
$wsST.Range("L2").Value2 = 45
$wsST.Range("M2").Value2 = 38

$wsST.Range("B3").Value2 = 57

# ---------------------------------------------------------------------
# TURNS sheet: update Road-row totals
# ---------------------------------------------------------------------
$wsTURNS = $wb.Worksheets.Item("TURNS")

$wsTURNS.Range("C3").Value2 = 6
$wsTURNS.Range("D3").Value2 = 7
$wsTURNS.Range("E3").Value2 = 6

# ---------------------------------------------------------------------
# PEN sheet: update penalty counts
# ---------------------------------------------------------------------
$wsPEN = $wb.Worksheets.Item("PEN")

$wsPEN.Range("B2").Value2 = 16
$wsPEN.Range("D2").Value2 = 15
$wsPEN.Range("B3").Value2 = 15
$wsPEN.Range("D4").Value2 = 9
